$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.188747525215149
$ws.Range("B1").Value = 2.169229745864868
$ws.Range("C1").Value = 3.744539737701416
$ws.Range("D1").Value = 3.280540704727173
$ws.Range("E1").Value = 1.139498829841614
